$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "Modelo" with same style as the other header cells (copy style from E1)
$ws.Range("F1").Value = "Modelo"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Update numeric values in row 2
$ws.Range("B2").Value = 0.08900161088493851
$ws.Range("C2").Value = 0.9987949731275695
$ws.Range("D2").Value = 0.2430300635827726

# Add new model info cell F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', AdaBoostRegressor(learning_rate=0.5))])"
